$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("L2").Value = "[4.697842502525113, 8.56626220691621]"
$ws.Range("M2").Value = [double]"4.657030316934652e-11"
$ws.Range("N2").Value = [double]"9.314060633869303e-11"
$ws.Range("P2").Value = "[-1.8365266363327724, -1.1824212590087715]"
$ws.Range("T2").Value = "[7.883681240561653, 10.295887844173162]"
$ws.Range("X2").Value = [double]"4.868428428428544"
$ws.Range("Y2").Value = [double]"7.56160160160178"

# Row 3
$ws.Range("L3").Value = "[2.9005384517494175, 10.330367437995653]"
$ws.Range("M3").Value = [double]"0.0005547531472123168"
$ws.Range("N3").Value = [double]"0.0005547531472123168"
$ws.Range("P3").Value = "[1.3899739268135027, 2.5723951858222724]"
$ws.Range("Q3").Value = [double]"3.811235771422616e-10"
$ws.Range("R3").Value = [double]"3.811235771422616e-10"
$ws.Range("T3").Value = "[7.473210708397952, 11.454946654825294]"
$ws.Range("X3").Value = [double]"13.2528528528529"
$ws.Range("Y3").Value = [double]"17.47579579579585"
